# Apply the changes described by the diff to "Route Cost RSO.xlsx"

$wb = $excel.ActiveWorkbook

# --- Sheet "Route" (sheet1) ---
$wsRoute = $wb.Worksheets.Item("Route")

# L3: date value -> text "14/4/2025" (new shared string)
$wsRoute.Range("L3").Value = "14/4/2025"

# D7: 130 -> 170 (I7/L7 formulas recalc automatically)
$wsRoute.Range("D7").Value = 170

# D9: 130 -> 200 (I9/L9 formulas recalc automatically)
$wsRoute.Range("D9").Value = 200

# D10: 200 -> 150 (I10/L10 formulas recalc automatically)
$wsRoute.Range("D10").Value = 150

# Update the selected cell to D7 (as in the saved workbook view)
$wsRoute.Range("D7").Select()

# --- Sheet "Mobil" (sheet2) ---
$wsMobil = $wb.Worksheets.Item("Mobil")

# G1: date 45689 -> 45748 (2025-04-01)
$wsMobil.Range("G1").Value = 45748
